{"js": "// Approved for publication, issue #166.\n// Insert a new \"Pre-Employment Screening and Vetting of External Candidates - FAQs\"\n// hyperlinked list item just before the \"User access\" heading (i.e. immediately\n// after the last item of the \"Personnel security clearances\" sub-list:\n// \"National Security Vetting questions\").\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph containing \"National Security Vetting questions\" - it is\n// the final entry of the numbered sub-list (ilvl=2) that the new entry must join.\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"National Security Vetting questions\") {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not locate anchor paragraph 'National Security Vetting questions'\");\n}\n\n// Inserting immediately after the anchor paragraph causes the new paragraph to\n// inherit the same list (numId/ilvl) and paragraph style (\"Compact\")\n// automatically - do NOT (re)assign .style/.styleBuiltIn afterwards, since doing\n// so strips the inherited list numbering that was just picked up.\nconst newParagraph = anchor.insertParagraph(\n  \"Pre-Employment Screening and Vetting of External Candidates - FAQs\",\n  Word.InsertLocation.after\n);\n\nawait context.sync();\n\n// Apply the hyperlink to the newly inserted paragraph's text.\nconst newRange = newParagraph.getRange();\nnewRange.hyperlink = \"pre-employment-screening-and-vetting-of-external-candidates-faqs.md\";\n\nawait context.sync();\n", "ps1": "# Approved for publication, issue #166.\n# Insert a new \"Pre-Employment Screening and Vetting of External Candidates - FAQs\"\n# hyperlinked list item just before the \"User access\" heading (i.e. immediately\n# after the last item of the \"Personnel security clearances\" sub-list:\n# \"National Security Vetting questions\").\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph: \"National Security Vetting questions\" - the last\n# entry of the numbered sub-list (ilvl=2, numId=1012) that the new entry must join.\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $ptext = $d.Paragraphs.Item($i).Range.Text\n    $ptext = $ptext.TrimEnd(\"`r\", \"`n\", [char]7)\n    if ($ptext -eq \"National Security Vetting questions\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate anchor paragraph 'National Security Vetting questions'\"\n}\n\n$anchor = $d.Paragraphs.Item($anchorIndex)\n$anchorRange = $anchor.Range\n$anchorRange.Collapse(0)          # wdCollapseEnd -> collapse to end of paragraph\n$anchorRange.InsertParagraphAfter()\n\n# The newly created (still empty) paragraph now sits right after the anchor and\n# has already inherited the same list numbering (ilvl=2, numId=1012) and the\n# \"Compact\" paragraph style, so there is no need to (re)apply them explicitly -\n# doing so via the Style property would actually strip the inherited numbering.\n$newParaIndex = $anchorIndex + 1\n$newPara = $d.Paragraphs.Item($newParaIndex)\n$newPara.Range.InsertBefore(\"Pre-Employment Screening and Vetting of External Candidates - FAQs\")\n\n# Build a fresh range over just the inserted text (excluding the trailing\n# paragraph mark) and turn it into a hyperlink.\n$newPara2 = $d.Paragraphs.Item($newParaIndex)\n$hyperlinkRange = $d.Range($newPara2.Range.Start, $newPara2.Range.End - 1)\n$d.Hyperlinks.Add($hyperlinkRange, \"pre-employment-screening-and-vetting-of-external-candidates-faqs.md\")\n"}
